$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.791318535804749
$ws.Range("B1").Value = 3.985836505889893
$ws.Range("C1").Value = 1.417524218559265
$ws.Range("D1").Value = 0.8544102907180786
$ws.Range("E1").Value = 0.4618767499923706
